$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1211
$ws1.Range("F4").Value = 2684
$ws1.Range("F5").Value = 240

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1211
$ws4.Range("F6").Value = 2684
$ws4.Range("F8").Value = 240
